$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: category changes from "Albatros" to "Prejuveniles";
# dia_1/dia_2 totals swap (F2 empties, G2 takes the 115 value).
$ws.Range("B2").Value = "Prejuveniles"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = 115

# Row 3: now holds the "caballeros" entry for Kern Pascuali (previously
# row 4's data), with dia_1 = 92 and dia_2 cleared.
$ws.Range("C3").Value = "caballeros"
$ws.Range("E3").Value = "Kern Pascuali, Juan Daniel"
$ws.Range("F3").Value = 92
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = 92

# Row 4 is removed entirely (its data was merged into row 3).
$ws.Rows(4).Delete()
